$d = $word.ActiveDocument

# Locate the paragraph containing the text we need to edit.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*diagramma UML dei casi*") {
        $target = $p
        break
    }
}

$pRange = $target.Range
$pStart = $pRange.Start

# The paragraph text is: "i diagramma UML dei casi d\u2019uso del software;"
# We need to change "diagramma" -> "diagrammi" by replacing the single
# letter "a" (right before " UML") with "i", while leaving the preceding
# "i " run intact and ending up with three runs:
#   1) " diagramm"                              (original run, truncated)
#   2) "i"                                       (new run)
#   3) " UML dei casi d\u2019uso del software;"          (new run)

$letterA = $d.Range($pStart + 10, $pStart + 11)
if ($letterA.Text -ne "a") {
    throw "Unexpected character at split point: [$($letterA.Text)]"
}

# Touching Font.Color (even re-applying the same resolved colour) forces
# the run to be carved out as its own element instead of being re-merged
# with its neighbours when the document is serialized back to OOXML.
$origColor = $letterA.Font.Color
$letterA.Font.Color = $origColor

# Replace the letter itself, turning "diagramma" into "diagrammi".
$letterA.Text = "i"

$r = $pRange
Write-Output $r.Text
